$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the old per-user email addresses with the new ones.
$ws.Range("G2").Value = "jack@gmail.com"
$ws.Range("G3").Value = "joe@gmail.com"
$ws.Range("G4").Value = "peter@gmail.com"

# Turn the email cells into mailto hyperlinks.
$ws.Hyperlinks.Add($ws.Range("G2"), "mailto:jack@gmail.com", [Type]::Missing, [Type]::Missing, "jack@gmail.com")
$ws.Hyperlinks.Add($ws.Range("G3"), "mailto:joe@gmail.com", [Type]::Missing, [Type]::Missing, "joe@gmail.com")
$ws.Hyperlinks.Add($ws.Range("G4"), "mailto:peter@gmail.com", [Type]::Missing, [Type]::Missing, "peter@gmail.com")

# Bump the header/data row heights slightly.
$ws.Rows("1:4").RowHeight = 16

# Move the active selection to G4.
$ws.Range("G4").Select() | Out-Null
